# Applies updated TPM-derived values to the LR-pairs sheet (Ptprz1-Ncam1),
# per refreshed NATMI script output ("update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.02270466666666667
$ws.Range("H2").Value = 0.06811400000000001
$ws.Range("I2").Value = 0.002206225855740089
$ws.Range("J2").Value = 0.002206225855740089
$ws.Range("M2").Value = 0.3912683333333333
$ws.Range("N2").Value = 1.173805
$ws.Range("O2").Value = 0.004923718964983145
$ws.Range("P2").Value = 0.004923718964983145
$ws.Range("Q2").Value = 0.008883617085555556
$ws.Range("R2").Value = 0.07995255377
$ws.Range("S2").Value = 0.00001086283608694364
$ws.Range("T2").Value = 0.00001086283608694364

# Row 3
$ws.Range("G3").Value = 0.02270466666666667
$ws.Range("H3").Value = 0.06811400000000001
$ws.Range("I3").Value = 0.002206225855740089
$ws.Range("J3").Value = 0.002206225855740089
$ws.Range("O3").Value = 0.07888477275715973
$ws.Range("P3").Value = 0.07888477275715973
$ws.Range("Q3").Value = 0.1423278054737778
$ws.Range("R3").Value = 1.280950249264
$ws.Range("S3").Value = 0.0001740376252810272
$ws.Range("T3").Value = 0.0001740376252810272

# Row 4
$ws.Range("G4").Value = 0.02270466666666667
$ws.Range("H4").Value = 0.06811400000000001
$ws.Range("I4").Value = 0.002206225855740089
$ws.Range("J4").Value = 0.002206225855740089
$ws.Range("M4").Value = 72.68848166666666
$ws.Range("N4").Value = 218.065445
$ws.Range("O4").Value = 0.9147115297293749
$ws.Range("P4").Value = 0.9147115297293749
$ws.Range("Q4").Value = 1.650367746747778
$ws.Range("R4").Value = 14.85330972073
$ws.Range("S4").Value = 0.002018060227432515
$ws.Range("T4").Value = 0.002018060227432515

# Row 5
$ws.Range("G5").Value = 0.02270466666666667
$ws.Range("H5").Value = 0.06811400000000001
$ws.Range("I5").Value = 0.002206225855740089
$ws.Range("J5").Value = 0.002206225855740089
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.117608
$ws.Range("N5").Value = 0.352824
$ws.Range("O5").Value = 0.001479978548482255
$ws.Range("P5").Value = 0.001479978548482255
$ws.Range("Q5").Value = 0.002670250437333333
$ws.Range("R5").Value = 0.024032253936
$ws.Range("S5").Value = 0.000003265166939602237
$ws.Range("T5").Value = 0.000003265166939602237

# Row 6
$ws.Range("I6").Value = 0.002281111990432972
$ws.Range("J6").Value = 0.002281111990432972
$ws.Range("M6").Value = 0.3912683333333333
$ws.Range("N6").Value = 1.173805
$ws.Range("O6").Value = 0.004923718964983145
$ws.Range("P6").Value = 0.004923718964983145
$ws.Range("Q6").Value = 0.009185154547777777
$ws.Range("R6").Value = 0.08266639093
$ws.Range("S6").Value = 0.00001123155436854528
$ws.Range("T6").Value = 0.00001123155436854528

# Row 7
$ws.Range("I7").Value = 0.002281111990432972
$ws.Range("J7").Value = 0.002281111990432972
$ws.Range("O7").Value = 0.07888477275715973
$ws.Range("P7").Value = 0.07888477275715973
$ws.Range("S7").Value = 0.0001799450009989374
$ws.Range("T7").Value = 0.0001799450009989374

# Row 8
$ws.Range("I8").Value = 0.002281111990432972
$ws.Range("J8").Value = 0.002281111990432972
$ws.Range("M8").Value = 72.68848166666666
$ws.Range("N8").Value = 218.065445
$ws.Range("O8").Value = 0.9147115297293749
$ws.Range("P8").Value = 0.9147115297293749
$ws.Range("Q8").Value = 1.706386336618889
$ws.Range("R8").Value = 15.35747702957
$ws.Range("S8").Value = 0.002086559438252963
$ws.Range("T8").Value = 0.002086559438252963

# Row 9
$ws.Range("I9").Value = 0.002281111990432972
$ws.Range("J9").Value = 0.002281111990432972
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.117608
$ws.Range("N9").Value = 0.352824
$ws.Range("O9").Value = 0.001479978548482255
$ws.Range("P9").Value = 0.001479978548482255
$ws.Range("Q9").Value = 0.002760887002666667
$ws.Range("R9").Value = 0.024847983024
$ws.Range("S9").Value = 0.000003375996812526458
$ws.Range("T9").Value = 0.000003375996812526458

# Row 10
$ws.Range("G10").Value = 10.24499966666667
$ws.Range("H10").Value = 30.734999
$ws.Range("I10").Value = 0.9955126621538269
$ws.Range("J10").Value = 0.9955126621538269
$ws.Range("M10").Value = 0.3912683333333333
$ws.Range("N10").Value = 1.173805
$ws.Range("O10").Value = 0.004923718964983145
$ws.Range("P10").Value = 0.004923718964983145
$ws.Range("Q10").Value = 4.008543944577222
$ws.Range("R10").Value = 36.076895501195
$ws.Range("S10").Value = 0.004901624574527657
$ws.Range("T10").Value = 0.004901624574527657

# Row 11
$ws.Range("G11").Value = 10.24499966666667
$ws.Range("H11").Value = 30.734999
$ws.Range("I11").Value = 0.9955126621538269
$ws.Range("J11").Value = 0.9955126621538269
$ws.Range("O11").Value = 0.07888477275715973
$ws.Range("P11").Value = 0.07888477275715973
$ws.Range("Q11").Value = 64.22240595044711
$ws.Range("R11").Value = 578.0016535540241
$ws.Range("S11").Value = 0.07853079013087977
$ws.Range("T11").Value = 0.07853079013087977

# Row 12
$ws.Range("G12").Value = 10.24499966666667
$ws.Range("H12").Value = 30.734999
$ws.Range("I12").Value = 0.9955126621538269
$ws.Range("J12").Value = 0.9955126621538269
$ws.Range("M12").Value = 72.68848166666666
$ws.Range("N12").Value = 218.065445
$ws.Range("O12").Value = 0.9147115297293749
$ws.Range("P12").Value = 0.9147115297293749
$ws.Range("Q12").Value = 744.693470445506
$ws.Range("R12").Value = 6702.241234009555
$ws.Range("S12").Value = 0.9106069100636894
$ws.Range("T12").Value = 0.9106069100636894

# Row 13
$ws.Range("G13").Value = 10.24499966666667
$ws.Range("H13").Value = 30.734999
$ws.Range("I13").Value = 0.9955126621538269
$ws.Range("J13").Value = 0.9955126621538269
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.117608
$ws.Range("N13").Value = 0.352824
$ws.Range("O13").Value = 0.001479978548482255
$ws.Range("P13").Value = 0.001479978548482255
$ws.Range("Q13").Value = 1.204893920797333
$ws.Range("R13").Value = 10.844045287176
$ws.Range("S13").Value = 0.001473337384730126
$ws.Range("T13").Value = 0.001473337384730126
